$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.70423055536827
$ws.Range("C2").Value = 0.711592290709234
$ws.Range("D2").Value = 0.06462893496711075
$ws.Range("E2").Value = 0.1267172724068644
$ws.Range("F2").Value = 2.986501631082831
$ws.Range("I2").Value = 1.363216454538488
$ws.Range("J2").Value = 0.2204158715208138
$ws.Range("B3").Value = 1.586601651986825
$ws.Range("C3").Value = 0.6620414870534432
$ws.Range("D3").Value = 0.06429010423621406
$ws.Range("E3").Value = 0.1246604846345072
$ws.Range("F3").Value = 2.938505611647656
$ws.Range("I3").Value = 1.350534565481766
$ws.Range("J3").Value = 0.2154193353751879
$ws.Range("B4").Value = 1.515334518320401
$ws.Range("C4").Value = 0.6320526742339894
$ws.Range("D4").Value = 0.06409995989472606
$ws.Range("E4").Value = 0.1234666764558909
$ws.Range("F4").Value = 2.910970408700294
$ws.Range("I4").Value = 1.343614453065236
$ws.Range("J4").Value = 0.2124921238354887
$ws.Range("B5").Value = 1.486532062299204
$ws.Range("C5").Value = 0.6199403482800676
$ws.Range("D5").Value = 0.06402702871551824
$ws.Range("E5").Value = 0.1229975236956768
$ws.Range("F5").Value = 2.900232230864745
$ws.Range("I5").Value = 1.341010528148978
$ws.Range("J5").Value = 0.2113344033914046
$ws.Range("B6").Value = 1.481763864135132
$ws.Range("C6").Value = 0.6179356167093601
$ws.Range("D6").Value = 0.06401519519099352
$ws.Range("E6").Value = 0.122920667204145
$ws.Range("F6").Value = 2.898478203933493
$ws.Range("I6").Value = 1.340591145240339
$ws.Range("J6").Value = 0.2111442805295667
$ws.Range("B7").Value = 1.514945110330757
$ws.Range("C7").Value = 0.6318888859562435
$ws.Range("D7").Value = 0.06409895780840102
$ws.Range("E7").Value = 0.1234602791568555
$ws.Range("F7").Value = 2.910823640692882
$ws.Range("I7").Value = 1.343578463156383
$ws.Range("J7").Value = 0.212476368429229
$ws.Range("B8").Value = 1.663472547840797
$ws.Range("C8").Value = 0.6944160393871357
$ws.Range("D8").Value = 0.06450842268699475
$ws.Range("E8").Value = 0.1259937350238651
$ws.Range("F8").Value = 2.969548761537453
$ws.Range("I8").Value = 1.358662755243714
$ws.Range("J8").Value = 0.2186637237224289
$ws.Range("B9").Value = 1.962409578619031
$ws.Range("C9").Value = 0.8205485658605198
$ws.Range("D9").Value = 0.06545117601579875
$ws.Range("E9").Value = 0.1315118951700534
$ws.Range("F9").Value = 3.10024043224584
$ws.Range("I9").Value = 1.395205310009516
$ws.Range("J9").Value = 0.2319248999538672
$ws.Range("B10").Value = 2.186854196268314
$ws.Range("C10").Value = 0.9154579630121589
$ws.Range("D10").Value = 0.06622621222108904
$ws.Range("E10").Value = 0.1359049820292952
$ws.Range("F10").Value = 3.20599492240612
$ws.Range("I10").Value = 1.42642111028627
$ws.Range("J10").Value = 0.242372866718668
$ws.Range("B11").Value = 2.290037890962537
$ws.Range("C11").Value = 0.9591429206124644
$ws.Range("D11").Value = 0.06659611178756109
$ws.Range("E11").Value = 0.1379779485215167
$ws.Range("F11").Value = 3.256278357772885
$ws.Range("I11").Value = 1.441597586395005
$ws.Range("J11").Value = 0.2472829179433376
$ws.Range("B12").Value = 2.3292687759008
$ws.Range("C12").Value = 0.9757602543112966
$ws.Range("D12").Value = 0.06673862651940254
$ws.Range("E12").Value = 0.1387737041175754
$ws.Range("F12").Value = 3.27563678475957
$ws.Range("I12").Value = 1.44748705910132
$ws.Range("J12").Value = 0.249165134036204
$ws.Range("B13").Value = 2.320812670831515
$ws.Range("C13").Value = 0.9721780648821436
$ws.Range("D13").Value = 0.06670782561716493
$ws.Range("E13").Value = 0.1386018441570016
$ws.Range("F13").Value = 3.271453423834856
$ws.Range("I13").Value = 1.446212285892898
$ws.Range("J13").Value = 0.2487587427599323
$ws.Range("B14").Value = 2.293262273343714
$ws.Range("C14").Value = 0.9605085314982489
$ws.Range("D14").Value = 0.06660778790896416
$ws.Range("E14").Value = 0.1380431997022171
$ws.Range("F14").Value = 3.257864605058927
$ws.Range("I14").Value = 1.4420792495077
$ws.Range("J14").Value = 0.2474373087047184
$ws.Range("B15").Value = 2.276407419069869
$ws.Range("C15").Value = 0.9533703878873325
$ws.Range("D15").Value = 0.06654682843422677
$ws.Range("E15").Value = 0.1377024177325907
$ws.Range("F15").Value = 3.249582498917164
$ws.Range("I15").Value = 1.439566261556749
$ws.Range("J15").Value = 0.2466308810141982
$ws.Range("B16").Value = 2.180132835181098
$ws.Range("C16").Value = 0.9126134354962119
$ws.Range("D16").Value = 0.066202382837389
$ws.Range("E16").Value = 0.1357710128124126
$ws.Range("F16").Value = 3.202752926076045
$ws.Range("I16").Value = 1.425449114453627
$ws.Range("J16").Value = 0.2420551706396452
$ws.Range("B17").Value = 2.121350035276521
$ws.Range("C17").Value = 0.8877420390746238
$ws.Range("D17").Value = 0.06599547948107443
$ws.Range("E17").Value = 0.1346052826549027
$ws.Range("F17").Value = 3.174584724031376
$ws.Range("I17").Value = 1.417040221467374
$ws.Range("J17").Value = 0.2392885843727868
$ws.Range("B18").Value = 2.087641522741251
$ws.Range("C18").Value = 0.8734846043910807
$ws.Range("D18").Value = 0.06587810844619213
$ws.Range("E18").Value = 0.1339417973468819
$ws.Range("F18").Value = 3.158587592075264
$ws.Range("I18").Value = 1.412295412065291
$ws.Range("J18").Value = 0.2377120983201308
$ws.Range("B19").Value = 2.076245831244421
$ws.Range("C19").Value = 0.8686654721012701
$ws.Range("D19").Value = 0.06583865087708318
$ws.Range("E19").Value = 0.1337183554010259
$ws.Range("F19").Value = 3.153206223576717
$ws.Range("I19").Value = 1.410704598339393
$ws.Range("J19").Value = 0.2371808575222047
$ws.Range("B20").Value = 2.127597016315462
$ws.Range("C20").Value = 0.8903846704398575
$ws.Range("D20").Value = 0.06601733592017922
$ws.Range("E20").Value = 0.1347286505595378
$ws.Range("F20").Value = 3.177562083442581
$ws.Range("I20").Value = 1.417925850269683
$ws.Range("J20").Value = 0.2395815604075153
$ws.Range("B21").Value = 2.301350205025415
$ws.Range("C21").Value = 0.9639341150630685
$ws.Range("D21").Value = 0.06663710553736735
$ws.Range("E21").Value = 0.1382069944567021
$ws.Range("F21").Value = 3.261847327021655
$ws.Range("I21").Value = 1.443289339129706
$ws.Range("J21").Value = 0.247824822546562
$ws.Range("B22").Value = 2.415826630459094
$ws.Range("C22").Value = 1.012439540589924
$ws.Range("D22").Value = 0.06705637066131231
$ws.Range("E22").Value = 0.1405430838728137
$ws.Range("F22").Value = 3.318783462534213
$ws.Range("I22").Value = 1.46069725699742
$ws.Range("J22").Value = 0.2533458026802435
$ws.Range("B23").Value = 2.354643743596228
$ws.Range("C23").Value = 0.9865108519380215
$ws.Range("D23").Value = 0.0668313172211441
$ws.Range("E23").Value = 0.1392905061520437
$ws.Range("F23").Value = 3.288224755190413
$ws.Range("I23").Value = 1.45132955417175
$ws.Range("J23").Value = 0.2503868400756772
$ws.Range("B24").Value = 2.124772486884865
$ws.Range("C24").Value = 0.8891898077601468
$ws.Range("D24").Value = 0.06600744969768613
$ws.Range("E24").Value = 0.1346728550381897
$ws.Range("F24").Value = 3.17621540552895
$ws.Range("I24").Value = 1.417525178639039
$ws.Range("J24").Value = 0.2394490621656189
$ws.Range("B25").Value = 1.880703683866045
$ws.Range("C25").Value = 0.7860400888342269
$ws.Range("D25").Value = 0.06518146600801344
$ws.Range("E25").Value = 0.1299598391042061
$ws.Range("F25").Value = 3.063194701740912
$ws.Range("I25").Value = 1.384561511340024
$ws.Range("J25").Value = 0.2282148771386829
